$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> C, D, E, F(text) values (A, B, G columns are unchanged)
$rows = @{
    2  = @(41, 122, 189, "[N/A,N/A]")
    3  = @(33, 86, 137, "[N/A,N/A]")
    4  = @(22, 158, 209, "[9,7]")
    5  = @(22, 94, 91, "[N/A,N/A]")
    6  = @(22, 110, 203, "[N/A,N/A]")
    7  = @(21, 178, 250, "[10,10]")
    8  = @(22, 90, 140, "[N/A,N/A]")
    9  = @(22, 137, 189, "[N/A,N/A]")
    10 = @(21, 141, 72, "[10,10]")
    11 = @(26, 107, 127, "[N/A,N/A]")
    12 = @(31, 86, 176, "[N/A,N/A]")
    13 = @(19, 139, 145, "[9,9]")
    14 = @(20, 78, 163, "[N/A,N/A]")
    15 = @(19, 70, 173, "[N/A,N/A]")
    16 = @(18, 116, 145, "[7,9]")
    17 = @(70, 304, 127, "[N/A,N/A]")
    18 = @(70, 322, 271, "[N/A,N/A]")
    19 = @(51, 476, 114, "[10,10]")
    20 = @(41, 155, 145, "[N/A,N/A]")
    21 = @(52, 146, 62, "[N/A,N/A]")
    22 = @(28, 195, 97, "[7,9]")
    23 = @(49, 245, 52, "[N/A,N/A]")
    24 = @(150, 476, 82, "[N/A,N/A]")
    25 = @(43, 290, 257, "[7,9]")
    26 = @(25, 26, 23, "[N/A,N/A]")
    27 = @(22, 62, 111, "[N/A,N/A]")
    28 = @(17, 58, 79, "[8,8]")
    29 = @(49, 213, 221, "[N/A,N/A]")
    30 = @(29, 86, 123, "[N/A,N/A]")
    31 = @(27, 154, 115, "[9,9]")
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 3).Value = $vals[0]
    $ws.Cells.Item($r, 4).Value = $vals[1]
    $ws.Cells.Item($r, 5).Value = $vals[2]
    $ws.Cells.Item($r, 6).Value = $vals[3]
}
